$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.405.27'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.311.30'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.24%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '186.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '578.22'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.96%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('E10').Value = '  +0.94%  '
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '3.887.53'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('E13').Value = '  -0.49%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.47'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '67.658.03'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E16').Value = '  -0.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.299.62'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '444.75'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.41%  '
$ws.Range('E19').Value = '  -0.91%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.95'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +3.49%  '
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.518'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.455.35'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E26').Value = '  +1.29%  '
$ws.Range('E27').Value = '  +0.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.06'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.35%  '
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '22.95'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.92%  '
$ws.Range('E32').Value = '  -2.40%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  -0.07%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.81'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.25%  '
$ws.Range('E36').Value = '  +4.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.78'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.72%  '
$ws.Range('E38').Value = '  -1.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '27.21'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('E40').Value = '  -0.98%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.759.10'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.53%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.14%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '24.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.68%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0674'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.42'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.21'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '326.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0274'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('E50').Value = '  +1.28%  '
$ws.Range('E51').Value = '  +1.57%  '
